$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08407699145275002
$ws.Range("H2").Value = 0.2495038727118119
$ws.Range("I2").Value = -10.87417598955746
$ws.Range("G3").Value = 0.09506055060865062
$ws.Range("H3").Value = -18.30178156260155
$ws.Range("G4").Value = -0.0803340366079255
$ws.Range("H4").Value = -420.7947151233894
$ws.Range("G5").Value = -0.03159704916454497
$ws.Range("H5").Value = 56.0246165583133
$ws.Range("G6").Value = 0.04167323674598129
$ws.Range("H6").Value = 17.23859064146115
$ws.Range("G7").Value = 0.02941654177588967
$ws.Range("H7").Value = 44.94407128134949
$ws.Range("G8").Value = -0.1539740972955305
$ws.Range("H8").Value = -9.125721111408348
$ws.Range("G9").Value = -0.1401472245702647
$ws.Range("H9").Value = -2.381682781198953
$ws.Range("G10").Value = -0.08264911212201734
$ws.Range("H10").Value = 23.49225849804611
$ws.Range("G11").Value = -0.08259246908123424
$ws.Range("H11").Value = -24.12293296052085
$ws.Range("G12").Value = -0.3870528409123202
$ws.Range("H12").Value = 6.778355422530064
$ws.Range("G13").Value = -0.5023355914908951
$ws.Range("H13").Value = -11.90382630924235
$ws.Range("G14").Value = -0.1013673296106214
$ws.Range("H14").Value = -99.79489309604989
$ws.Range("G15").Value = 0.03268301130920041
$ws.Range("H15").Value = 139.5240575490354
$ws.Range("G16").Value = 0.1284386701063236
$ws.Range("H16").Value = -12.04197525420225
$ws.Range("G17").Value = 0.1518004176716578
$ws.Range("H17").Value = 23.79593781628978
$ws.Range("G18").Value = 0.1207208176661175
$ws.Range("H18").Value = -12.80999865492479
$ws.Range("G19").Value = 0.1317670490046652
$ws.Range("H19").Value = 38.11083196311385
$ws.Range("G20").Value = 0.04108054488050385
$ws.Range("H20").Value = 60.2577097412198
$ws.Range("G21").Value = 0.03634350327225206
$ws.Range("H21").Value = -51.64972811824406
$ws.Range("G24").Value = 0.1112481793794625
$ws.Range("H24").Value = 10.76468381168089
$ws.Range("G25").Value = 0.1401739454478299
$ws.Range("H25").Value = -7.5043541116513
$ws.Range("G26").Value = 0.06757116481109801
$ws.Range("H26").Value = -14.60169201293565
$ws.Range("G27").Value = 0.07280490676445721
$ws.Range("H27").Value = -27.0876357783589
$ws.Range("G28").Value = -0.2511424702665349
$ws.Range("H28").Value = -17.8413246874822
$ws.Range("G29").Value = -0.2221356457444656
$ws.Range("H29").Value = -8.196991075493266
$ws.Range("G30").Value = 0.05440179504616999
$ws.Range("H30").Value = 23.27162954057643
$ws.Range("G31").Value = 0.0202199629916218
$ws.Range("H31").Value = -23.2146845965622
$ws.Range("G32").Value = 0.1147697996080134
$ws.Range("H32").Value = 20.86055677440321
$ws.Range("G33").Value = 0.1148658987707932
$ws.Range("H33").Value = 10.48528222782648
$ws.Range("G34").Value = 0.006739757988254588
$ws.Range("H34").Value = -85.48332299268742
$ws.Range("G35").Value = 0.006582357565931249
$ws.Range("H35").Value = -13.12171008397028
$ws.Range("G36").Value = 0.05956327806959027
$ws.Range("H36").Value = 3.164445680930076
$ws.Range("G37").Value = 0.01330912456872372
$ws.Range("H37").Value = -81.07514374755516
$ws.Range("G38").Value = 0.00593796508307031
$ws.Range("H38").Value = -88.66506613459941
$ws.Range("G39").Value = 0.04721442479372522
$ws.Range("H39").Value = 127.67842498533
$ws.Range("G40").Value = -0.001755392655731389
$ws.Range("H40").Value = 79.32640652271277
$ws.Range("G41").Value = 0.03052246486384081
$ws.Range("H41").Value = -13.67068985437523
$ws.Range("G42").Value = 0.1289266584298991
$ws.Range("H42").Value = -3.560180480433381
$ws.Range("G43").Value = 0.1494987561138102
$ws.Range("H43").Value = 0.3506742362397638
$ws.Range("G44").Value = -0.001494510762795267
$ws.Range("H44").Value = 82.44123772150789
$ws.Range("G45").Value = 0.00404352435017862
$ws.Range("H45").Value = 136.8267217143798
$ws.Range("G46").Value = -0.02747829481516197
$ws.Range("H46").Value = -734.434003561243
$ws.Range("G47").Value = -0.005188237060141277
$ws.Range("H47").Value = 44.08312085906672
$ws.Range("G48").Value = 0.06605263008404766
$ws.Range("H48").Value = 31.39020311637638
$ws.Range("G49").Value = 0.06842467974841773
$ws.Range("H49").Value = 3.572093913333144
$ws.Range("G50").Value = 0.1427981168115234
$ws.Range("H50").Value = -11.44724609088292
$ws.Range("G51").Value = 0.1400383794659987
$ws.Range("H51").Value = -18.16122921202611
$ws.Range("G52").Value = -0.1752230513571305
$ws.Range("H52").Value = -9.220194096309964
$ws.Range("G53").Value = -0.144913957043227
$ws.Range("H53").Value = -14.96135134884795
$ws.Range("G54").Value = 0.09126263002600364
$ws.Range("H54").Value = -2.624522311904939
$ws.Range("G55").Value = 0.1208028971734723
$ws.Range("H55").Value = 6.828881803034531
$ws.Range("G56").Value = -0.0073564917752215
$ws.Range("H56").Value = -0.7632107213147814
$ws.Range("G57").Value = -0.02343193090678332
$ws.Range("H57").Value = -2.468754129558146
$ws.Range("G58").Value = 0.03185962551019064
$ws.Range("H58").Value = -43.49925336520885
$ws.Range("G59").Value = 0.07433733047692882
$ws.Range("H59").Value = 3.505161656407051
$ws.Range("G60").Value = 0.06637762220727771
$ws.Range("H60").Value = -5.139528438934647
$ws.Range("G61").Value = 0.0791779447804694
$ws.Range("H61").Value = 66.59646533457367
$ws.Range("G62").Value = 0.07224223508638043
$ws.Range("H62").Value = -0.9842231387238702
$ws.Range("G63").Value = 0.08046696407190683
$ws.Range("H63").Value = 23.05039259599913
$ws.Range("G64").Value = -0.03535014524168426
$ws.Range("H64").Value = 14.64514785049825
$ws.Range("G65").Value = 0.01803228992548703
$ws.Range("H65").Value = 136.5590722390017
$ws.Range("G66").Value = 0.03136840711406265
$ws.Range("H66").Value = 65.67258823980045
$ws.Range("G67").Value = 0.05162462178504321
$ws.Range("H67").Value = 97.41815345711319
$ws.Range("G68").Value = -0.001515185505012106
$ws.Range("H68").Value = -365.8437633921408
$ws.Range("G69").Value = -0.01414314016496307
$ws.Range("H69").Value = -9.532250063827426
$ws.Range("G70").Value = -0.0482052701084016
$ws.Range("H70").Value = -75.62354444538688
$ws.Range("G71").Value = -0.06732241194170754
$ws.Range("H71").Value = -22.20008682244654
$ws.Range("G72").Value = -0.1623879228114489
$ws.Range("H72").Value = -9.47440423920645
$ws.Range("G73").Value = -0.1501816591960114
$ws.Range("H73").Value = -3.723896194444685
$ws.Range("G74").Value = 0.1466494249418579
$ws.Range("H74").Value = 16.34124445170325
$ws.Range("G75").Value = 0.1429515578838513
$ws.Range("H75").Value = 5.75532061338556
$ws.Range("G76").Value = -0.0854431989084769
$ws.Range("H76").Value = -148.1077160862452
$ws.Range("G77").Value = -0.09863998301548565
$ws.Range("H77").Value = -113.5665149770866
$ws.Range("G78").Value = 0.09570286514399941
$ws.Range("H78").Value = 3.83125532801944
$ws.Range("G79").Value = 0.07576976962448762
$ws.Range("H79").Value = -21.48570792186607
$ws.Range("G80").Value = -0.1570960312002562
$ws.Range("H80").Value = 3.273649170886771
$ws.Range("G81").Value = -0.1677007349458716
$ws.Range("H81").Value = 22.52163102415944
$ws.Range("G82").Value = 0.1475813847152427
$ws.Range("H82").Value = 6.367814464665305
$ws.Range("G83").Value = 0.2040110508529877
$ws.Range("H83").Value = 23.928777056746
$ws.Range("G84").Value = 0.04946741685997214
$ws.Range("H84").Value = 253.464233261784
$ws.Range("G85").Value = 0.03036770679534123
$ws.Range("H85").Value = 34.13781694556874
